# Update ChanjoKe FHIR IG - StructureDefinition-batch-number.xlsx
$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# URL
$wsMeta.Range("B2").Value = "https://intellisoft-consulting.github.io/igs/ChanjoKe-FHIR-IG/StructureDefinition/batch-number"

# Title (shared between Metadata!B5 and Elements!L2)
$wsMeta.Range("B5").Value = "Extension for Batch Number "
$wsElem.Range("L2").Value = "Extension for Batch Number "

# Date
$wsMeta.Range("B8").Value = "2024-08-27T20:30:12+00:00"

# Publisher
$wsMeta.Range("B9").Value = "Intellisoft Consulting Ltd"

# Contact
$wsMeta.Range("B10").Value = "Intellisoft Consulting Ltd (https://www.intellisoftkenya.com/, info[at]intellisoftkenya.com)"

# Description (shared between Metadata!B12 and Elements!M2)
$wsMeta.Range("B12").Value = "AExtension for Batch Number "
$wsElem.Range("M2").Value = "AExtension for Batch Number "

# Context
$wsMeta.Range("B21").Value = "element:SupplyDelivery"
